# Apply the "Eco model" results update to the F1 Logistics Roundtrip sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8: visual separator row (all three cells = "---", same style as header rows) ---
$ws.Range("A1").Copy()
$ws.Range("A8:C8").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("A8").Value = "---"
$ws.Range("B8").Value = "---"
$ws.Range("C8").Value = "---"

# --- Row 9 / Row 10: new Eco_NN / Eco_IH result rows ---
# Column A cells reuse the same label formatting as the rest of column A.
$ws.Range("A1").Copy()
$ws.Range("A9").PasteSpecial(-4122)
$ws.Range("A1").Copy()
$ws.Range("A10").PasteSpecial(-4122)

# Fill in the values in the same order the original author typed them, so the
# shared-string table comes out in the same sequence.
$ws.Range("A10").Value = "Eco_IH"
$ws.Range("A9").Value = "Eco_NN"
$ws.Range("B9").Value = "78,147 km"
$ws.Range("B10").Value = "73,793 km"
$ws.Range("C9").Value = "65,270 km"
$ws.Range("C10").Value = "65,714 km"

# Columns B/C on the new rows use a distinct (but visually identical, plain black Arial 10)
# explicit-color font, matching the source workbook's style.
$dataCells = @("B9", "C9", "B10", "C10")
foreach ($addr in $dataCells) {
    $rng = $ws.Range($addr)
    $rng.Font.Color = 0
    $rng.Font.Name = "Arial"
    $rng.Font.Size = 10
}
